# Regen - now constrained
#
# Adds a new "data_collection_mode list" validation-source sheet (DDA / DIA),
# positioned right after "ion_mobility list" (pushing column_length_unit list
# and everything after it down by one slot), and wires a new list-type data
# validation onto column X ("data_collection_mode") of the "Export as TSV"
# sheet, matching the other "must come from list" validations already on
# that sheet.

$wb = $excel.ActiveWorkbook

# 1) Insert the new lookup-list worksheet right after "ion_mobility list".
$afterSheet = $wb.Worksheets.Item("ion_mobility list")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "data_collection_mode list"
$newSheet.Range("A1").Value = "DDA"
$newSheet.Range("A2").Value = "DIA"

# 2) Add the matching data validation to column X on the main sheet.
$ws = $wb.Worksheets.Item("Export as TSV")
$range = $ws.Range("X2:X1048576")
$range.Validation.Add(3, 1, 1, "='data_collection_mode list'!`$A`$1:`$A`$2")
$range.Validation.ErrorTitle = "Value must come from list"
$range.Validation.ErrorMessage = "Value must be one of: DDA / DIA."
$range.Validation.InputTitle = ""
$range.Validation.InputMessage = ""
$range.Validation.ShowInput = $true
$range.Validation.ShowError = $true
$range.Validation.IgnoreBlank = $true

# Restore the originally-active sheet/selection (adding the new sheet made it active).
$ws.Activate()
[void]$ws.Range("A1").Select()
